# repull data, push all data, mean calculation
# Update column F ("dSF") values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -9
$ws.Range("F3").Value  = -9
$ws.Range("F5").Value  = -7
$ws.Range("F6").Value  = -2
$ws.Range("F7").Value  = -4
$ws.Range("F9").Value  = -7
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = -3
